## Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3.
##
## The upstream change merely re-serializes every run-level bold toggle
## (`<w:b w:val="true"/>`) to the spelling POI 5.2.3 emits for the same
## boolean value (`<w:b w:val="on"/>`) -- "true" and "on" are both valid
## ST_OnOff spellings for True, so the runs stay (and must stay) bold.
## Re-apply Bold explicitly to every one of those red "Invalid block: ..."
## runs so the property round-trips through Word's object model.

$d = $word.ActiveDocument

$searchText = "Invalid block: Unexpected tag EOF missing [ENDFOR] while parsing m:for v | self.eClassifiers"

$rng = $d.Content
$rng.Find.ClearFormatting()

while ($rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Font.Bold = $true
    $rng.Collapse(0)
}
